# "completed a few more AIs, added a few new ones as well"
#
# This script:
#  1. Marks row 9 ("install latest version of python (activestate
#     python)") as Done by setting B9 to "Yes", matching the formatting
#     already used by the other "Yes" cells in column B.
#  2. Fills in row 13 with a brand-new AI task: Order 7, "upload latest
#     stock predictions", owner "OA".
#  3. Leaves the selection on A13 (where the new row was entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: mark "install latest version of python" as Done -------------
# Copy the formatting already used by the other checked "Yes" cells in
# column B (e.g. B3) onto B9 before writing the value, so the new cell
# matches the rest of the column instead of keeping its old, unused style.
$ws.Range("B3").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Yes"

# --- Row 13: new AI task --------------------------------------------------
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "upload latest stock predictions"
$ws.Range("E13").Value = "OA"

# --- Selection -------------------------------------------------------------
$ws.Range("A13").Select()
